$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-5 from 45183 (2023-09-14)
# to 45184 (2023-09-15), keeping existing number formatting.
$ws.Range("C2:C5").Value = 45184
